# Project Spyn Gantt Chart - "fix mistakes" update
#
# Semantic changes being applied:
#   1. Project start date (Q1) is switched from a volatile =TODAY() formula
#      to a fixed literal date (10/3/2023, serial 45202).
#   2. Two duplicate/leftover "Milestone 2" task rows (15 & 16) are
#      corrected to use the right task text, matching the later
#      "Milestone 3" tasks, with updated progress % and end-date formulas.
#   3. Progress percentages are tweaked on rows 17-19.
#   4. End-date formulas are adjusted on rows 18 & 19.
#   5. Rows 21-25 start dates are normalized to the same formula
#      ($E$9+15) instead of a mix of relative formula / hard-coded dates.
#   6. The sheet's scroll position / selection is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project schedule")

# --- 1. Project start date: was =TODAY(), now a fixed date -------------
$ws.Range("Q1").Value = 45202

# --- 2. Rows 15 & 16: correct the task text + numbers -------------------
$ws.Range("B15").Value = "Make sure the Car Doesn't Hurt Guy"
$ws.Range("D15").Value = 1
$ws.Range("F15").Formula = "=E15"

$ws.Range("B16").Value = "Figure out how to Navigate The Car Around the Maze"
$ws.Range("D16").Value = 0.7
$ws.Range("F16").Formula = "=E16+3"

# --- 3. Rows 17-19: progress % + end-date formula tweaks ---------------
$ws.Range("D17").Value = 0.85

$ws.Range("D18").Value = 0.9
$ws.Range("F18").Formula = "=E18+2"

$ws.Range("D19").Value = 1
$ws.Range("E19").Formula = "=E18-2"
$ws.Range("F19").Formula = "=E19"

# --- 4. Rows 21-25: normalize start-date formulas -----------------------
$ws.Range("E21").Formula = "=`$E`$9+15"
$ws.Range("E22").Formula = "=`$E`$9+15"
$ws.Range("E23").Formula = "=`$E`$9+15"
$ws.Range("E24").Formula = "=`$E`$9+15"
$ws.Range("E25").Formula = "=`$E`$9+15"

# --- 5. Sheet scroll position / selection --------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("E28").Select()
